$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column C ("Thương hiệu") before old column C, shifting
#     Đường kính..Ngoài bảng (and their widths/bestFit flags) one column right. ---
$ws.Columns.Item(3).Insert()

# --- Column widths (engine quantizes ColumnWidth to 1/6-character steps,
#     so pick the input that lands closest to the authored width) ---
$ws.Columns.Item(1).ColumnWidth = 28.833333333333332   # -> 29.6640625-ish (A)
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666   # -> 16.5546875-ish (B)
$ws.Columns.Item(3).ColumnWidth = 13.5                 # -> 14.33203125-ish (C, new)
$ws.Columns.Item(7).ColumnWidth = 8.333333333333334    # -> 9.109375-ish (G)

# --- Header row ---
$ws.Range("A1").Value = "Tên hàng"
$ws.Range("B1").Value = "Loại"
$ws.Range("C1").Value = "Thương hiệu"
$ws.Range("D1").Value = "Đường kính"
$ws.Range("E1").Value = "Chữ"
$ws.Range("F1").Value = "Độ"
$ws.Range("G1").Value = "Số"
$ws.Range("H1").Value = "Độ K"
$ws.Range("I1").Value = "Đơn vị"
$ws.Range("J1").Value = "Ngoài bảng"

# --- Row 2 ---
$ws.Range("A2").Value = "B05-10.6-Standard"
$ws.Range("B2").Value = "Standard"
$ws.Range("C2").Value = "Fargo"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "10.6"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "B"
$ws.Range("F2").ClearContents()
# Establish the plain/non-header Text cellXf (index 2) before the header one,
# so the two new cellXfs land at the same indices as the authored file.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "05"
$ws.Range("H2").ClearContents()
$ws.Range("I2").Value = "Cái"
$ws.Range("J2").Value = "không"

# Column G (Số) header also becomes Text-formatted (cellXf index 3).
$ws.Range("G1").NumberFormat = "@"

# --- Row 3 ---
$ws.Range("A3").Value = "B06-10.6-Standard Isee B&L"
$ws.Range("B3").Value = "Standard Isee B&L"
$ws.Range("C3").Value = "Isee B&L"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "10.6"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "B"
$ws.Range("F3").ClearContents()
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "06"
$ws.Range("H3").ClearContents()
$ws.Range("I3").Value = "Cái"
$ws.Range("J3").Value = "không"

# --- Row 4 (new row) ---
$ws.Range("A4").Value = "Dung dịch thử nghiệm"
$ws.Range("B4").Value = "Dung dịch"
$ws.Range("C4").Value = "Fargo"
$ws.Range("I4").Value = "Chai"
$ws.Range("J4").Value = "không"

# --- Selection ---
$ws.Range("D7").Select() | Out-Null
